$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(465, 44539, 6, 16, 303.030303030303),
    @(466, 44540, 18, 33, 625),
    @(467, 44541, 1, 30, 568.1818181818182),
    @(468, 44542, 7, 35, 662.8787878787879),
    @(469, 44543, 8, 42, 795.4545454545454),
    @(470, 44544, 9, 50, 946.969696969697),
    @(471, 44545, 0, 49, 928.0303030303031),
    @(472, 44546, 10, 53, 1003.787878787879),
    @(473, 44547, 2, 37, 700.7575757575758),
    @(474, 44548, 0, 36, 681.8181818181818),
    @(475, 44550, 4, 33, 625),
    @(476, 44551, 4, 29, 549.2424242424242),
    @(477, 44552, 1, 21, 397.7272727272727),
    @(478, 44553, 6, 27, 511.3636363636364),
    @(479, 44554, 7, 24, 454.5454545454545),
    @(480, 44555, 4, 26, 492.4242424242424),
    @(481, 44556, 5, 31, 587.1212121212121),
    @(482, 44557, 4, 31, 587.1212121212121),
    @(483, 44558, 0, 27, 511.3636363636364),
    @(484, 44559, 1, 27, 511.3636363636364),
    @(485, 44560, 7, 28, 530.3030303030304),
    @(486, 44561, 11, 32, 606.060606060606),
    @(487, 44562, 11, 39, 738.6363636363636),
    @(488, 44563, 3, 37, 700.7575757575758),
    @(489, 44564, 6, 39, 738.6363636363636),
    @(490, 44565, 2, 41, 776.5151515151515),
    @(491, 44566, 8, 48, 909.090909090909)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Copy the date-column style (s="2") from row 464 (A column) down to the new A-column cells
$ws.Range("A464").Copy() | Out-Null
$ws.Range("A465:A491").PasteSpecial(-4122) | Out-Null

# Restore the values that PasteSpecial may have touched (formats only, so values remain, but re-set to be safe)
$ws.Cells.Item(465, 1).Value = 44539
$ws.Cells.Item(466, 1).Value = 44540
$ws.Cells.Item(467, 1).Value = 44541
$ws.Cells.Item(468, 1).Value = 44542
$ws.Cells.Item(469, 1).Value = 44543
$ws.Cells.Item(470, 1).Value = 44544
$ws.Cells.Item(471, 1).Value = 44545
$ws.Cells.Item(472, 1).Value = 44546
$ws.Cells.Item(473, 1).Value = 44547
$ws.Cells.Item(474, 1).Value = 44548
$ws.Cells.Item(475, 1).Value = 44550
$ws.Cells.Item(476, 1).Value = 44551
$ws.Cells.Item(477, 1).Value = 44552
$ws.Cells.Item(478, 1).Value = 44553
$ws.Cells.Item(479, 1).Value = 44554
$ws.Cells.Item(480, 1).Value = 44555
$ws.Cells.Item(481, 1).Value = 44556
$ws.Cells.Item(482, 1).Value = 44557
$ws.Cells.Item(483, 1).Value = 44558
$ws.Cells.Item(484, 1).Value = 44559
$ws.Cells.Item(485, 1).Value = 44560
$ws.Cells.Item(486, 1).Value = 44561
$ws.Cells.Item(487, 1).Value = 44562
$ws.Cells.Item(488, 1).Value = 44563
$ws.Cells.Item(489, 1).Value = 44564
$ws.Cells.Item(490, 1).Value = 44565
$ws.Cells.Item(491, 1).Value = 44566
